$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 2-21 (B and C columns) with new mean/std values ---
$existing = @(
    @(0.48359375, 0.02335924435582196),
    @(0.58359375, 0.05438532588736598),
    @(0.86015625, 0.04670542243010976),
    @(0.8875, 0.07995512657500456),
    @(0.84140625, 0.08851255648917841),
    @(0.5984375, 0.07499186153760286),
    @(0.78828125, 0.07903377689712798),
    @(0.66015625, 0.08999647345521379),
    @(0.96640625, 0.02217979619777873),
    @(0.54140625, 0.0832901906551726),
    @(0.4437654613049464, 0.03481862921325921),
    @(0.545300820285353, 0.0914016193680357),
    @(0.8578398247430734, 0.04821188574262342),
    @(0.8851801621745971, 0.0825056703533834),
    @(0.8223808028758519, 0.1186395578679433),
    @(0.5577213654494683, 0.09726491516147848),
    @(0.7809473994724301, 0.08461696668921222),
    @(0.6218493503981242, 0.1168062037608873),
    @(0.9662073127486972, 0.02254374849569267),
    @(0.5138071992584632, 0.1007392482368634)
)

for ($i = 0; $i -lt $existing.Count; $i++) {
    $row = 2 + $i
    $ws.Cells.Item($row, 2).Value = $existing[$i][0]
    $ws.Cells.Item($row, 3).Value = $existing[$i][1]
}

# --- Append new rows 22-31 (index 20-29) ---
# First, replicate the formatting of row 21's A cell (bold, bordered, centered)
# down through the new rows in column A.
$ws.Range("A21").Copy() | Out-Null
$ws.Range("A22:A31").PasteSpecial(-4122) | Out-Null

$newRows = @(
    @(20, 0.50653076171875, 0.03120027400532249),
    @(21, 0.61748046875, 0.07355619537380011),
    @(22, 0.9482421875, 0.02229199103241224),
    @(23, 0.94569091796875, 0.05747054755514734),
    @(24, 0.97191162109375, 0.01485858978539397),
    @(25, 0.67701416015625, 0.05131886461119931),
    @(26, 0.8916748046875, 0.07269655975796875),
    @(27, 0.7451416015625, 0.1018909249993786),
    @(28, 0.99129638671875, 0.01170485626575681),
    @(29, 0.52017822265625, 0.08253417116575971)
)

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $row = 22 + $i
    $ws.Cells.Item($row, 1).Value = $newRows[$i][0]
    $ws.Cells.Item($row, 2).Value = $newRows[$i][1]
    $ws.Cells.Item($row, 3).Value = $newRows[$i][2]
}
